# LIVEHTA-1269 testdata correction
# Rebuild the PRISMA sample data table with the new column layout/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table (previously A1:I6) completely, then drop the two
# columns (H, I) that are no longer part of the new, narrower table.
$ws.Range("A1:I6").Clear()
$ws.Range("H1:I1").EntireColumn.Delete()

# ---- Header row (except F1, added later - see below) ------------------
$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Prisma_Population"
$ws.Cells.Item(1, 3).Value = "Prisma_Excel_File"
$ws.Cells.Item(1, 4).Value = "Study_Types"
$ws.Cells.Item(1, 5).Value = "stdy_type_locators"
$ws.Cells.Item(1, 7).Value = "Prisma_Image"

# ---- Data rows ----------------------------------------------------------
$data = @(
    @{ Row = 2;  A = "pop1"; B = "LIVEHTA Automation - Test_NonOncology_Automation_3"; C = "\Testdata\Templates\PRISMA\Test_Sachin\12. PRISMA_Pfizer_IC AML Mylotarg.xlsx"; D = "Clinical-Interventional";    E = "total_record_number";                F = 100;  G = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Clinical.PNG" },
    @{ Row = 3;  A = "pop1";                                                                                                                              D = "Clinical-RWE";               E = "total_excluded_number";              F = 200;  G = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Clinical.PNG" },
    @{ Row = 4;  A = "pop1";                                                                                                                              D = "Economic";                   E = "total_screenedTiAb_number";          F = 300;  G = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Econ.PNG" },
    @{ Row = 5;  A = "pop1";                                                                                                                              D = "Quality of life";            E = "total_excluded_screenedTiAb_number"; F = 400;  G = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_QoL.PNG" },
    @{ Row = 6;  A = "pop1";                                                                                                                              D = "Real-world Evidence";        E = "fulltext_review";                    F = 500;  G = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Clinical.PNG" },
    @{ Row = 7;  A = "pop1";                                                                                                                                                                E = "excluded_fulltext_review";           F = 600 },
    @{ Row = 8;  A = "pop1";                                                                                                                                                                E = "fulltext_exclusion_reason";          F = 700 },
    @{ Row = 9;  A = "pop1";                                                                                                                                                                E = "total_greyliterature_number";        F = 800 },
    @{ Row = 10; A = "pop1";                                                                                                                                                                E = "original_studies";                   F = 900 },
    @{ Row = 11; A = "pop1";                                                                                                                                                                E = "records_number";                     F = 1000 }
)

# Column E (stdy_type_locators values) is filled top-to-bottom first...
foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 5).Value = $entry.E
}

# ...then the remaining columns A-D and G are populated...
foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    if ($entry.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $entry.B }
    if ($entry.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $entry.C }
    if ($entry.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $entry.D }
    if ($entry.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $entry.F }
    if ($entry.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $entry.G }
}

# ...and finally the F1 header ("stdy_type_values") is added last.
$ws.Cells.Item(1, 6).Value = "stdy_type_values"

# ---- Header row styling (re-apply the centered style used before) --------
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").VerticalAlignment = -4108

# ---- Column widths --------------------------------------------------------
# Columns A-D keep their original widths (unchanged by this edit), so they
# are intentionally left alone. Only E, F and G need new widths.
$ws.Columns.Item(5).ColumnWidth = 32.77734375
$ws.Columns.Item(6).ColumnWidth = 20.5546875
$ws.Columns.Item(7).ColumnWidth = 71.6640625

# ---- Selection -------------------------------------------------------------
$ws.Range("E16").Select()
